$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Order number
# ------------------------------------------------------------------
$d.Content.Find.Execute("No.  124", $true, $false, $false, $false, $false, $true, 1, $false, "No.  129", 2)

# ------------------------------------------------------------------
# 2. The original "Fecha devoluci\u00f3n" value (2019-02-28) must become
#    2019-03-01 BEFORE the "Fecha"/"Fecha evento" (2019-02-27 ->
#    2019-02-28) replacement below, otherwise it would get caught by
#    that later ReplaceAll too.
# ------------------------------------------------------------------
$d.Content.Find.Execute("2019-02-28", $true, $false, $false, $false, $false, $true, 1, $false, "2019-03-01", 2)

# ------------------------------------------------------------------
# 3. Fecha / Fecha evento: both occurrences change identically, so
#    ReplaceAll is safe now.
# ------------------------------------------------------------------
$d.Content.Find.Execute("2019-02-27", $true, $false, $false, $false, $false, $true, 1, $false, "2019-02-28", 2)

# ------------------------------------------------------------------
# 4. Client name
# ------------------------------------------------------------------
$d.Content.Find.Execute("Pepito Perez", $true, $false, $false, $false, $false, $true, 1, $false, "Diego Gallardo", 2)

# ------------------------------------------------------------------
# 5. Address
# ------------------------------------------------------------------
$d.Content.Find.Execute("Barrio Centenario", $true, $false, $false, $false, $false, $true, 1, $false, "cra 1 # 23-4", 2)

# ------------------------------------------------------------------
# 6. Phone (Tel)
# ------------------------------------------------------------------
$d.Content.Find.Execute("3203259689", $true, $false, $false, $false, $false, $true, 1, $false, "320556657", 2)

# ------------------------------------------------------------------
# 7. Referencia (case fix nadie -> Nadie)
# ------------------------------------------------------------------
$d.Content.Find.Execute("nadie", $true, $false, $false, $false, $false, $true, 1, $false, "Nadie", 2)

# ------------------------------------------------------------------
# 8. Celular
# ------------------------------------------------------------------
$d.Content.Find.Execute("34342234234", $true, $false, $false, $false, $false, $true, 1, $false, "320254544", 2)

# ------------------------------------------------------------------
# 9. Concepto line
# ------------------------------------------------------------------
$d.Content.Find.Execute("Pantalon 3 - (p4)", $true, $false, $false, $false, $false, $true, 1, $false, "(p5) - Producto 5", 2)

# ------------------------------------------------------------------
# 10. Abono: "29997" -> "0". In the source document this run sits next
#     to two other runs (" - " and "CANCELADO") that share the exact
#     same run formatting (sz 20). This engine coalesces adjacent runs
#     with identical formatting whenever a paragraph is edited, which
#     would wrongly fold " - CANCELADO" into the edited run. To keep
#     them as independent runs (matching the original structure) we
#     momentarily flip Bold on/off on the neighbouring text immediately
#     after the text edit -- toggling a character property forces the
#     engine to keep/re-establish a run boundary there, and clearing it
#     straight back to its original (falsy) state leaves no visible
#     trace in the resulting rPr.
# ------------------------------------------------------------------
$d.Content.Find.Execute("29997", $true, $false, $false, $false, $false, $true, 1, $false, "0", 2)

$found = $d.Content
$found.Find.Execute("0 - CANCELADO")
$zeroStart = $found.Start
$zeroEnd = $zeroStart + 1

$splitA = $d.Range($zeroStart, $zeroEnd)
$splitA.Bold = 1
$splitA2 = $d.Range($zeroStart, $zeroEnd)
$splitA2.Bold = 0

$cancel = $d.Content
$cancel.Find.Execute("CANCELADO")
$splitB = $d.Range($cancel.Start, $cancel.End)
$splitB.Bold = 1
$splitB2 = $d.Range($cancel.Start, $cancel.End)
$splitB2.Bold = 0

# ------------------------------------------------------------------
# 11. Fecha devoluci\u00f3n row: Saldo value 20000 -> 42500
# ------------------------------------------------------------------
$d.Content.Find.Execute("20000", $true, $false, $false, $false, $false, $true, 1, $false, "42500", 2)

# ------------------------------------------------------------------
# 12. TOTAL now shows the discount percentage applied
# ------------------------------------------------------------------
$d.Content.Find.Execute("49997", $true, $false, $false, $false, $false, $true, 1, $false, "42500 ( %15)", 2)
